$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The batsman name in the source data ends with a non-breaking space
# (U+00A0), matching the existing rows 2:4 exactly.
$nbsp = [char]0x00A0
$batsman = "Mohammed Shami" + $nbsp

# New rows to append (5,6,7) - duplicate the existing match data.
$rows = @(
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Kings XI Punjab", "Mumbai Indians", $batsman, "2", "2", "0", "0", "100.00"),
    @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Kings XI Punjab", "Delhi Capitals", $batsman, "0", "0", "0", "0", "-"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Kings XI Punjab", "Sunrisers Hyderabad", $batsman, "0", "1", "0", "0", "0.00")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    # Force the numeric-looking columns (G:K) to be stored as text, matching
    # the rest of the sheet (all cells are text-typed, per numberStoredAsText).
    $ws.Range("A" + $r + ":K" + $r).NumberFormat = "@"

    for ($col = 0; $col -lt $rowData.Length; $col++) {
        $colLetter = [string]([char](65 + $col))
        $cellRef = $colLetter + $r
        $ws.Range($cellRef).Value = $rowData[$col]
    }
}
